$wb = $excel.ActiveWorkbook

# --- Sheet "Typography" ---
$typo = $wb.Worksheets.Item("Typography")

# Row 5: Wildcard Characters (H) was mistakenly filled; move value to
# Wildcard Ranges (I) instead, matching the new convention used by the
# other rows below.
$typo.Range("H5").Value = $null
$typo.Range("I5").Value = "0x20-0x7E"

# Rows 6 and 7: fill in the Wildcard Ranges column (I) the same way.
$typo.Range("I6").Value = "0x20-0x7E"
$typo.Range("I7").Value = "0x20-0x7E"

# Row 8: brand new Typography entry "Typography_01" (RPM gauge / timer font).
$typo.Range("B8").Value = "Typography_01"
$typo.Range("C8").Value = "verdana.ttf"
$typo.Range("D8").Value = 26
$typo.Range("E8").Value = 4
$typo.Range("F8").Value = "?"
$typo.Range("I8").Value = "0x20-0x7E"
# Columns B:E inherit a non-default column style; reset the newly written
# cells back to the workbook default so the row matches the plain (no
# explicit cell style) look of the other data rows above it.
$typo.Range("B8:E8").Style = "Normal"

# --- Sheet "Translation" ---
$trans = $wb.Worksheets.Item("Translation")

# Row 5: RPM text now includes unit suffix "rpm(s)".
$trans.Range("F5").Value = "RPM: <value> rpm(s)"

# Row 6: Input Voltage text now includes unit suffix "V".
$trans.Range("F6").Value = "Input Voltage: <value> V"

# Row 7: repurposed as the "Timer" text using the new Typography_00 style.
$trans.Range("C7").Value = "Typography_00"
$trans.Range("F7").Value = "Timer"

# Row 8: brand new translation entry for the Typography_01 style.
$trans.Range("B8").Value = "SingleUseId5"
$trans.Range("C8").Value = "Typography_01"
$trans.Range("D8").Value = "Left"
$trans.Range("E8").Value = "LTR"
$trans.Range("F8").Value = "A<value>"
